$wb = $excel.ActiveWorkbook

# 1. Insert new "Player Info" sheet at the very front of the workbook
$firstSheet = $wb.Worksheets.Item(1)
$piSheet = $wb.Worksheets.Add($firstSheet)
$piSheet.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $piSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$piSheet.Range("A1:D1").Font.Bold = $true
$piSheet.Range("A1:D1").HorizontalAlignment = -4108
$piSheet.Range("A1:D1").VerticalAlignment = -4160
$piSheet.Range("A1:D1").Borders.LineStyle = 1

$piSheet.Range("A2").NumberFormat = "@"
$piSheet.Range("A2").Value = "4855"
$piSheet.Range("B2").Value = "Shane Charles Getkate"
$piSheet.Range("C2").Value = "Right Handed"
$piSheet.Range("D2").Value = "Right Arm Medium Fast"

# 2. Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2:D5").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4343"
$battingSheet.Range("D3").Value = "4347"
$battingSheet.Range("D4").Value = "4352"
$battingSheet.Range("D5").Value = "4496"

# 3. Update "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2:B5").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4343"
$bowlingSheet.Range("B3").Value = "4347"
$bowlingSheet.Range("B4").Value = "4352"
$bowlingSheet.Range("B5").Value = "4496"

Write-Output "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
